# Update the "dSF" (column F) values for the rows that were repulled/recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = -3
    5  = 4
    6  = -2
    8  = 5
    10 = -7
    13 = -2
    15 = 1
    18 = 0
    21 = -4
    22 = -3
    24 = 5
    25 = 0
    28 = 2
    34 = -1
    36 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
